$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @"
2,46048,455,625,1080,1,26.01.20261
3,46048.01041666666,379,579,958,2,26.01.20262
4,46048.02083333334,381,578,959,3,26.01.20263
5,46048.03125,379,0,379,4,26.01.20264
6,46048.04166666666,382,564,946,5,26.01.20265
7,46048.05208333334,380,563,943,6,26.01.20266
8,46048.0625,383,564,947,7,26.01.20267
9,46048.07291666666,380,563,943,8,26.01.20268
10,46048.08333333334,435,560,995,9,26.01.20269
11,46048.09375,433,558,991,10,26.01.202610
12,46048.10416666666,436,559,995,11,26.01.202611
13,46048.11458333334,434,558,992,12,26.01.202612
14,46048.125,422,553,975,13,26.01.202613
15,46048.13541666666,383,530,913,14,26.01.202614
16,46048.14583333334,382,529,911,15,26.01.202615
17,46048.15625,385,0,385,16,26.01.202616
18,46048.16666666666,391,0,391,17,26.01.202617
19,46048.17708333334,393,536,929,18,26.01.202618
20,46048.1875,391,535,926,19,26.01.202619
21,46048.19791666666,394,537,931,20,26.01.202620
22,46048.20833333334,393,573,966,21,26.01.202621
23,46048.21875,396,590,986,22,26.01.202622
24,46048.22916666666,397,596,993,23,26.01.202623
25,46048.23958333334,423,615,1038,24,26.01.202624
26,46048.25,692,1024,1716,25,26.01.202625
27,46048.26041666666,755,1076,1831,26,26.01.202626
28,46048.27083333334,834,1123,1957,27,26.01.202627
29,46048.28125,832,1134,1966,28,26.01.202628
30,46048.29166666666,923,1260,2183,29,26.01.202629
31,46048.30208333334,785,1174,1959,30,26.01.202630
32,46048.3125,0,1159,1159,31,26.01.202631
33,46048.32291666666,786,0,786,32,26.01.202632
34,46048.33333333334,0,1169,1169,33,26.01.202633
35,46048.34375,785,1166,1951,34,26.01.202634
36,46048.35416666666,788,0,788,35,26.01.202635
37,46048.36458333334,0,1167,1167,36,26.01.202636
38,46048.375,747,1137,1884,37,26.01.202637
39,46048.38541666666,814,1193,2007,38,26.01.202638
40,46048.39583333334,824,1210,2034,39,26.01.202639
41,46048.40625,887,1249,2136,40,26.01.202640
42,46048.41666666666,871,1236,2107,41,26.01.202641
43,46048.42708333334,858,1234,2092,42,26.01.202642
44,46048.4375,0,0,0,43,26.01.202643
45,46048.44791666666,0,0,0,44,26.01.202644
46,46048.45833333334,0,0,0,45,26.01.202645
47,46048.46875,0,0,0,46,26.01.202646
48,46048.47916666666,0,0,0,47,26.01.202647
49,46048.48958333334,0,0,0,48,26.01.202648
50,46048.5,0,0,0,49,26.01.202649
51,46048.51041666666,0,0,0,50,26.01.202650
52,46048.52083333334,0,0,0,51,26.01.202651
53,46048.53125,0,0,0,52,26.01.202652
54,46048.54166666666,0,0,0,53,26.01.202653
55,46048.55208333334,0,0,0,54,26.01.202654
56,46048.5625,0,0,0,55,26.01.202655
57,46048.57291666666,0,0,0,56,26.01.202656
58,46048.58333333334,0,0,0,57,26.01.202657
59,46048.59375,0,0,0,58,26.01.202658
60,46048.60416666666,0,0,0,59,26.01.202659
61,46048.61458333334,0,0,0,60,26.01.202660
62,46048.625,0,0,0,61,26.01.202661
63,46048.63541666666,0,0,0,62,26.01.202662
64,46048.64583333334,0,0,0,63,26.01.202663
65,46048.65625,0,0,0,64,26.01.202664
66,46048.66666666666,0,0,0,65,26.01.202665
67,46048.67708333334,0,0,0,66,26.01.202666
68,46048.6875,0,0,0,67,26.01.202667
69,46048.69791666666,0,0,0,68,26.01.202668
70,46048.70833333334,0,0,0,69,26.01.202669
71,46048.71875,0,0,0,70,26.01.202670
72,46048.72916666666,0,0,0,71,26.01.202671
73,46048.73958333334,0,0,0,72,26.01.202672
74,46048.75,0,0,0,73,26.01.202673
75,46048.76041666666,0,0,0,74,26.01.202674
76,46048.77083333334,0,0,0,75,26.01.202675
77,46048.78125,0,0,0,76,26.01.202676
78,46048.79166666666,0,0,0,77,26.01.202677
79,46048.80208333334,0,0,0,78,26.01.202678
80,46048.8125,0,0,0,79,26.01.202679
81,46048.82291666666,0,0,0,80,26.01.202680
82,46048.83333333334,0,0,0,81,26.01.202681
83,46048.84375,0,0,0,82,26.01.202682
84,46048.85416666666,0,0,0,83,26.01.202683
85,46048.86458333334,0,0,0,84,26.01.202684
86,46048.875,0,0,0,85,26.01.202685
87,46048.88541666666,0,0,0,86,26.01.202686
88,46048.89583333334,0,0,0,87,26.01.202687
89,46048.90625,0,0,0,88,26.01.202688
90,46048.91666666666,0,0,0,89,26.01.202689
91,46048.92708333334,0,0,0,90,26.01.202690
92,46048.9375,0,0,0,91,26.01.202691
93,46048.94791666666,0,0,0,92,26.01.202692
94,46048.95833333334,0,0,0,93,26.01.202693
95,46048.96875,0,0,0,94,26.01.202694
96,46048.97916666666,0,0,0,95,26.01.202695
97,46048.98958333334,0,0,0,96,26.01.202696
98,46049,369,593,962,1,27.01.20261
99,46049.01041666666,365,522,887,2,27.01.20262
100,46049.02083333334,350,603,953,3,27.01.20263
101,46049.03125,344,616,960,4,27.01.20264
102,46049.04166666666,368,584,952,5,27.01.20265
103,46049.05208333334,366,581,947,6,27.01.20266
104,46049.0625,368,586,954,7,27.01.20267
105,46049.07291666666,365,584,949,8,27.01.20268
106,46049.08333333334,370,586,956,9,27.01.20269
107,46049.09375,368,585,953,10,27.01.202610
108,46049.10416666666,300,472,772,11,27.01.202611
109,46049.11458333334,274,438,712,12,27.01.202612
110,46049.125,252,450,702,13,27.01.202613
111,46049.13541666666,244,452,696,14,27.01.202614
112,46049.14583333334,248,588,836,15,27.01.202615
113,46049.15625,251,665,916,16,27.01.202616
114,46049.16666666666,257,723,980,17,27.01.202617
115,46049.17708333334,255,732,987,18,27.01.202618
116,46049.1875,0,735,735,19,27.01.202619
117,46049.19791666666,263,0,263,20,27.01.202620
118,46049.20833333334,342,691,1033,21,27.01.202621
119,46049.21875,375,690,1065,22,27.01.202622
120,46049.22916666666,392,685,1077,23,27.01.202623
121,46049.23958333334,399,694,1093,24,27.01.202624
122,46049.25,403,1016,1419,25,27.01.202625
123,46049.26041666666,396,1034,1430,26,27.01.202626
124,46049.27083333334,399,1043,1442,27,27.01.202627
125,46049.28125,442,1045,1487,28,27.01.202628
126,46049.29166666666,760,1209,1969,29,27.01.202629
127,46049.30208333334,770,1225,1995,30,27.01.202630
128,46049.3125,767,1220,1987,31,27.01.202631
129,46049.32291666666,755,0,755,32,27.01.202632
130,46049.33333333334,789,1224,2013,33,27.01.202633
131,46049.34375,791,1220,2011,34,27.01.202634
132,46049.35416666666,0,1242,1242,35,27.01.202635
133,46049.36458333334,785,1250,2035,36,27.01.202636
134,46049.375,637,1358,1995,37,27.01.202637
135,46049.38541666666,636,1371,2007,38,27.01.202638
136,46049.39583333334,635,1373,2008,39,27.01.202639
137,46049.40625,633,1371,2004,40,27.01.202640
138,46049.41666666666,600,1335,1935,41,27.01.202641
139,46049.42708333334,599,0,599,42,27.01.202642
140,46049.4375,0,1333,1333,43,27.01.202643
141,46049.44791666666,0,0,0,44,27.01.202644
142,46049.45833333334,0,0,0,45,27.01.202645
143,46049.46875,0,0,0,46,27.01.202646
144,46049.47916666666,0,0,0,47,27.01.202647
145,46049.48958333334,0,0,0,48,27.01.202648
146,46049.5,0,0,0,49,27.01.202649
147,46049.51041666666,0,0,0,50,27.01.202650
148,46049.52083333334,0,0,0,51,27.01.202651
149,46049.53125,0,0,0,52,27.01.202652
150,46049.54166666666,0,0,0,53,27.01.202653
151,46049.55208333334,0,0,0,54,27.01.202654
152,46049.5625,0,0,0,55,27.01.202655
153,46049.57291666666,0,0,0,56,27.01.202656
154,46049.58333333334,0,0,0,57,27.01.202657
155,46049.59375,0,0,0,58,27.01.202658
156,46049.60416666666,0,0,0,59,27.01.202659
157,46049.61458333334,0,0,0,60,27.01.202660
158,46049.625,0,0,0,61,27.01.202661
159,46049.63541666666,0,0,0,62,27.01.202662
160,46049.64583333334,0,0,0,63,27.01.202663
161,46049.65625,0,0,0,64,27.01.202664
162,46049.66666666666,0,0,0,65,27.01.202665
163,46049.67708333334,0,0,0,66,27.01.202666
164,46049.6875,0,0,0,67,27.01.202667
165,46049.69791666666,0,0,0,68,27.01.202668
166,46049.70833333334,0,0,0,69,27.01.202669
167,46049.71875,0,0,0,70,27.01.202670
168,46049.72916666666,0,0,0,71,27.01.202671
169,46049.73958333334,0,0,0,72,27.01.202672
170,46049.75,0,0,0,73,27.01.202673
171,46049.76041666666,0,0,0,74,27.01.202674
172,46049.77083333334,0,0,0,75,27.01.202675
173,46049.78125,0,0,0,76,27.01.202676
174,46049.79166666666,0,0,0,77,27.01.202677
175,46049.80208333334,0,0,0,78,27.01.202678
176,46049.8125,0,0,0,79,27.01.202679
177,46049.82291666666,0,0,0,80,27.01.202680
178,46049.83333333334,0,0,0,81,27.01.202681
179,46049.84375,0,0,0,82,27.01.202682
180,46049.85416666666,0,0,0,83,27.01.202683
181,46049.86458333334,0,0,0,84,27.01.202684
182,46049.875,0,0,0,85,27.01.202685
183,46049.88541666666,0,0,0,86,27.01.202686
184,46049.89583333334,0,0,0,87,27.01.202687
185,46049.90625,0,0,0,88,27.01.202688
186,46049.91666666666,0,0,0,89,27.01.202689
187,46049.92708333334,0,0,0,90,27.01.202690
188,46049.9375,0,0,0,91,27.01.202691
189,46049.94791666666,0,0,0,92,27.01.202692
190,46049.95833333334,0,0,0,93,27.01.202693
191,46049.96875,0,0,0,94,27.01.202694
192,46049.97916666666,0,0,0,95,27.01.202695
193,46049.98958333334,0,0,0,96,27.01.202696
"@

$lines = $data -split "`n" | Where-Object { $_.Trim() -ne "" }

$numRows = $lines.Count
$arrA = New-Object 'object[,]' $numRows,1
$arrBCDE = New-Object 'object[,]' $numRows,4
$arrF = New-Object 'object[,]' $numRows,1

$startRow = 2
$endRow = $startRow + $numRows - 1

for ($i = 0; $i -lt $numRows; $i++) {
    $parts = $lines[$i].Trim() -split ","
    $rowNum = [int]$parts[0]
    $aVal = [double]$parts[1]
    $bVal = [double]$parts[2]
    $cVal = [double]$parts[3]
    $dVal = [double]$parts[4]
    $eVal = [double]$parts[5]
    $fVal = $parts[6]

    $arrA[$i,0] = $aVal
    $arrBCDE[$i,0] = $bVal
    $arrBCDE[$i,1] = $cVal
    $arrBCDE[$i,2] = $dVal
    $arrBCDE[$i,3] = $eVal
    $arrF[$i,0] = $fVal
}

$ws.Range("A$startRow`:A$endRow").Value2 = $arrA
$ws.Range("B$startRow`:E$endRow").Value2 = $arrBCDE
$ws.Range("F$startRow`:F$endRow").Value = $arrF
